$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data row (row 2) : strip the old text-number-format styling from A2 first ---
$ws.Range("A2").Style = "Normal"
$ws.Range("A2").Value = 1771176464

# --- Header row (row 1) ---
$ws.Range("A1").Value = "Phone Numbers"
$ws.Range("B1").Value = "Message"
$ws.Range("C1").Value = "Message Status"

# Header styling: solid yellow fill across A1:C1
$ws.Range("A1:C1").Interior.Color = 65535

# --- Remaining data row values ---
$ws.Range("B2").Value = "Sent"
$ws.Range("C2").Value = "not seen"

# --- Column widths (approximate target sizing for the new 3-column layout) ---
$ws.Columns.Item(1).ColumnWidth = 18.5
$ws.Columns.Item(2).ColumnWidth = 17.666666666666668
$ws.Columns.Item(3).ColumnWidth = 17.166666666666668
